# Natmi LR-pair result sheet (Lgi3-Adam23): refresh the sending/target-cluster
# cross table following Dr Hou's advice -- the ligand/receptor-expressing-cell
# counts and every derived statistic (detection rate, expression, specificity,
# edge weights, ...) change, and the table grows from 3 data rows to 8 data rows
# (FAPs/sCs as sending clusters x ECs/FAPs/M2/sCs as target clusters).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row layout: A=Sending cluster, B=Ligand symbol, C=Receptor symbol, D=Target cluster,
# E..T = the 16 numeric NATMI statistics columns (same order as the header row).
# Row 2: FAPs -> ECs
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Lgi3"
$ws.Range("C2").Value = "Adam23"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.6167776666666667
$ws.Range("H2").Value = 1.850333
$ws.Range("I2").Value = 0.7836323164322263
$ws.Range("J2").Value = 0.7836323164322262
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.269164
$ws.Range("N2").Value = 0.8074920000000001
$ws.Range("O2").Value = 0.01089095165781685
$ws.Range("P2").Value = 0.01089095165781686
$ws.Range("Q2").Value = 0.1660143438706667
$ws.Range("R2").Value = 1.494129094836
$ws.Range("S2").Value = 0.008534501675766417
$ws.Range("T2").Value = 0.008534501675766417

# Row 3: FAPs -> FAPs
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Lgi3"
$ws.Range("C3").Value = "Adam23"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.6167776666666667
$ws.Range("H3").Value = 1.850333
$ws.Range("I3").Value = 0.7836323164322263
$ws.Range("J3").Value = 0.7836323164322262
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 20.07911633333334
$ws.Range("N3").Value = 60.237349
$ws.Range("O3").Value = 0.8124440315867432
$ws.Range("P3").Value = 0.8124440315867433
$ws.Range("Q3").Value = 12.38435052080189
$ws.Range("R3").Value = 111.459154687217
$ws.Range("S3").Value = 0.6366573984438564
$ws.Range("T3").Value = 0.6366573984438564

# Row 4: FAPs -> M2
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Lgi3"
$ws.Range("C4").Value = "Adam23"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.6167776666666667
$ws.Range("H4").Value = 1.850333
$ws.Range("I4").Value = 0.7836323164322263
$ws.Range("J4").Value = 0.7836323164322262
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.174593
$ws.Range("N4").Value = 0.523779
$ws.Range("O4").Value = 0.007064406543197522
$ws.Range("P4").Value = 0.007064406543197523
$ws.Range("Q4").Value = 0.1076850631563333
$ws.Range("R4").Value = 0.969165568407
$ws.Range("S4").Value = 0.005535897263664851
$ws.Range("T4").Value = 0.005535897263664851

# Row 5: FAPs -> sCs
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Lgi3"
$ws.Range("C5").Value = "Adam23"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.6167776666666667
$ws.Range("H5").Value = 1.850333
$ws.Range("I5").Value = 0.7836323164322263
$ws.Range("J5").Value = 0.7836323164322262
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 4.191587666666666
$ws.Range("N5").Value = 12.574763
$ws.Range("O5").Value = 0.1696006102122423
$ws.Range("P5").Value = 0.1696006102122424
$ws.Range("Q5").Value = 2.585277660675444
$ws.Range("R5").Value = 23.267498946079
$ws.Range("S5").Value = 0.1329045190489386
$ws.Range("T5").Value = 0.1329045190489386

# Row 6: sCs -> ECs
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Lgi3"
$ws.Range("C6").Value = "Adam23"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.1702976666666667
$ws.Range("H6").Value = 0.510893
$ws.Range("I6").Value = 0.2163676835677737
$ws.Range("J6").Value = 0.2163676835677737
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.269164
$ws.Range("N6").Value = 0.8074920000000001
$ws.Range("O6").Value = 0.01089095165781685
$ws.Range("P6").Value = 0.01089095165781686
$ws.Range("Q6").Value = 0.04583800115066667
$ws.Range("R6").Value = 0.4125420103560001
$ws.Range("S6").Value = 0.002356449982050438
$ws.Range("T6").Value = 0.002356449982050438

# Row 7: sCs -> FAPs
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Lgi3"
$ws.Range("C7").Value = "Adam23"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.1702976666666667
$ws.Range("H7").Value = 0.510893
$ws.Range("I7").Value = 0.2163676835677737
$ws.Range("J7").Value = 0.2163676835677737
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 20.07911633333334
$ws.Range("N7").Value = 60.237349
$ws.Range("O7").Value = 0.8124440315867432
$ws.Range("P7").Value = 0.8124440315867433
$ws.Range("Q7").Value = 3.419426660295223
$ws.Range("R7").Value = 30.774839942657
$ws.Range("S7").Value = 0.1757866331428868
$ws.Range("T7").Value = 0.1757866331428868

# Row 8: sCs -> M2
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Lgi3"
$ws.Range("C8").Value = "Adam23"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.1702976666666667
$ws.Range("H8").Value = 0.510893
$ws.Range("I8").Value = 0.2163676835677737
$ws.Range("J8").Value = 0.2163676835677737
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.174593
$ws.Range("N8").Value = 0.523779
$ws.Range("O8").Value = 0.007064406543197522
$ws.Range("P8").Value = 0.007064406543197523
$ws.Range("Q8").Value = 0.02973278051633333
$ws.Range("R8").Value = 0.267595024647
$ws.Range("S8").Value = 0.001528509279532671
$ws.Range("T8").Value = 0.001528509279532672

# Row 9: sCs -> sCs
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Lgi3"
$ws.Range("C9").Value = "Adam23"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.1702976666666667
$ws.Range("H9").Value = 0.510893
$ws.Range("I9").Value = 0.2163676835677737
$ws.Range("J9").Value = 0.2163676835677737
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 4.191587666666666
$ws.Range("N9").Value = 12.574763
$ws.Range("O9").Value = 0.1696006102122423
$ws.Range("P9").Value = 0.1696006102122424
$ws.Range("Q9").Value = 0.7138175992621111
$ws.Range("R9").Value = 6.424358393359
$ws.Range("S9").Value = 0.03669609116330378
$ws.Range("T9").Value = 0.03669609116330379
